# "30/11/2017 MAMATHA CHICK IN"
#
# 1) The original "SUN Nov 26 ... 2017" timestamp paragraph was split across
#    two runs (an artifact of how it was typed/pasted in originally). Word
#    naturally coalesces a run when you search-and-replace its exact text,
#    so doing that collapses the two runs into the single run seen in the
#    target document.
# 2) A brand new purchase entry ("TUE Nov 28 10:25:47 PST 2017" / MANJU /
#    BEET / cash) is appended to the bottom of the document, mirroring the
#    layout/formatting of the existing entry above it (same PlainText /
#    Courier New styling, bold banner line, red "Amount Received" line).

$d = $word.ActiveDocument

# --- 1. Merge the two "SUN Nov 26" / " 10:37:50 PST 2017" runs into one ---
$d.Content.Find.Execute(
    "SUN Nov 26 10:37:50 PST 2017", $false, $false, $false, $false, $false,
    $true, 1, $false, "SUN Nov 26 10:37:50 PST 2017", 2) | Out-Null

# --- 2. Append the new "TUE Nov 28" purchase block -------------------------
# Anchor on the existing bold "Amount balance ... - 636.0" paragraph (the
# last paragraph that has real content) and insert the new paragraphs right
# after it, ahead of the trailing blank bold paragraph that already closes
# the document.
$amountBalance = $d.Paragraphs(10)
$insertionPoint = $d.Range($amountBalance.Range.End, $amountBalance.Range.End)

$courierRPr = '<w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>'
$redRPr = '<w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="FF0000"/>'

$newParagraphsXml = @"
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$courierRPr<w:b/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$courierRPr</w:rPr></w:pPr>
  <w:r><w:rPr>$courierRPr</w:rPr><w:t>TUE Nov 28</w:t></w:r>
  <w:r><w:rPr>$courierRPr</w:rPr><w:t xml:space="preserve"> 10:25:47 PST 2017</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$courierRPr</w:rPr></w:pPr>
  <w:r><w:rPr>$courierRPr</w:rPr><w:t>Person Name</w:t></w:r>
  <w:r><w:rPr>$courierRPr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$courierRPr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$courierRPr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$courierRPr</w:rPr><w:tab/><w:t>- MANJU</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$courierRPr</w:rPr></w:pPr>
  <w:r><w:rPr>$courierRPr</w:rPr><w:t>---------------------------------------------------------------</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$courierRPr</w:rPr></w:pPr>
  <w:r><w:rPr>$courierRPr</w:rPr><w:t>Item Name</w:t></w:r>
  <w:r><w:rPr>$courierRPr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$courierRPr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$courierRPr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$courierRPr</w:rPr><w:tab/><w:t>- BEET</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$redRPr</w:rPr></w:pPr>
  <w:r><w:rPr>$redRPr</w:rPr><w:t>Amount Received</w:t></w:r>
  <w:r><w:rPr>$redRPr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$redRPr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$redRPr</w:rPr><w:tab/><w:t>- 636</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$courierRPr</w:rPr></w:pPr>
  <w:r><w:rPr>$courierRPr</w:rPr><w:t>Amount Received mode</w:t></w:r>
  <w:r><w:rPr>$courierRPr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$courierRPr</w:rPr><w:tab/><w:t>- CASH AND CLEARD</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$courierRPr</w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$courierRPr<w:b/></w:rPr></w:pPr></w:p>
"@

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $newParagraphsXml + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($packageXml) | Out-Null

Write-Host "Inserted TUE Nov 28 purchase entry; merged SUN Nov 26 timestamp run."
